# Generate Report for Handoff
# Adds a new file's localization-status row (3ea21f8f-5947-46ae-af78-6e626f390fe3.md)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet (table "Overview")
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A4").Value = "3ea21f8f-5947-46ae-af78-6e626f390fe3.md"
$wsOverview.Range("B4").Value = "e2e\3ea21f8f-5947-46ae-af78-6e626f390fe3.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2017-02-09 06:20:36"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$hlOverview = $wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a944f47c61b96bc1acc4df6404c91a4503624b2c/e2e/3ea21f8f-5947-46ae-af78-6e626f390fe3.md", "", "", "e2e\3ea21f8f-5947-46ae-af78-6e626f390fe3.md")

# ---------------------------------------------------------------------------
# zh-cn sheet (table "zh-cn")
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$rowZhCn = $loZhCn.ListRows.Add()

$wsZhCn.Range("A4").Value = "3ea21f8f-5947-46ae-af78-6e626f390fe3.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "3ea21f8f-5947-46ae-af78-6e626f390fe3.73fd09fc2679f07d38a056a5802eb3cc15a6b4e4.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2017-02-09 06:20:18"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = ""
$wsZhCn.Range("L4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M4").Value = ""
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "True"
$wsZhCn.Range("P4").Value = ""
$wsZhCn.Range("Q4").Value = "False"
$wsZhCn.Range("R4").Value = ""

$hlZhCn = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a944f47c61b96bc1acc4df6404c91a4503624b2c/e2e/3ea21f8f-5947-46ae-af78-6e626f390fe3.md", "", "", "3ea21f8f-5947-46ae-af78-6e626f390fe3.md")

# ---------------------------------------------------------------------------
# de-de sheet (table "de-de")
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$rowDeDe = $loDeDe.ListRows.Add()

$wsDeDe.Range("A4").Value = "3ea21f8f-5947-46ae-af78-6e626f390fe3.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "3ea21f8f-5947-46ae-af78-6e626f390fe3.73fd09fc2679f07d38a056a5802eb3cc15a6b4e4.de-de.xlf"
$wsDeDe.Range("H4").Value = "2017-02-09 06:20:36"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = ""
$wsDeDe.Range("L4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M4").Value = ""
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "True"
$wsDeDe.Range("P4").Value = ""
$wsDeDe.Range("Q4").Value = "False"
$wsDeDe.Range("R4").Value = ""

$hlDeDe = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a944f47c61b96bc1acc4df6404c91a4503624b2c/e2e/3ea21f8f-5947-46ae-af78-6e626f390fe3.md", "", "", "3ea21f8f-5947-46ae-af78-6e626f390fe3.md")

# ---------------------------------------------------------------------------
# Re-apply the HyperLink cell style last, after the hyperlink objects are
# created, so the blue-underline look used elsewhere in the sheets is kept.
# ---------------------------------------------------------------------------
$wsOverview.Range("B4").Style = "HyperLink"
$wsZhCn.Range("A4").Style = "HyperLink"
$wsDeDe.Range("A4").Style = "HyperLink"
